$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Seed the new shared strings in the same order they were first introduced
# (gbp/jpy, then the two "opt len and correl" labels, then "Re-date").
$ws.Range("B16").Value = "gbp/jpy"
$ws.Range("A16").Value = "opt len and correl1"
$ws.Range("A18").Value = "opt len and correl2"
$ws.Range("A17").Value = "Re-date"

# --- Row 16 : gbp/jpy, opt len and correl1, H1, D/E, L/M (re-dated), N/O/P ---
$ws.Range("C16").Value = "H1"
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 93
$ws.Range("L16").Value = 42614
$ws.Range("L16").NumberFormat = "mmm-yy"
$ws.Range("M16").Value = 42887
$ws.Range("M16").NumberFormat = "mmm-yy"
$ws.Range("N16").Value = 52
$ws.Range("O16").Value = 2600
$ws.Range("P16").Value = 1.74

# --- Row 17 : Re-date ---
$ws.Range("L17").Value = 41518
$ws.Range("L17").NumberFormat = "mmm-yy"
$ws.Range("M17").Value = 42156
$ws.Range("M17").NumberFormat = "mmm-yy"
$ws.Range("N17").Value = 63
$ws.Range("O17").Value = -1000
$ws.Range("P17").Value = 0.76

# --- Row 18 : gbp/jpy, opt len and correl2, H1, D/E/G, L/M (re-dated), N/O/P ---
$ws.Range("B18").Value = "gbp/jpy"
$ws.Range("C18").Value = "H1"
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = 93
$ws.Range("G18").Value = 35
$ws.Range("L18").Value = 42614
$ws.Range("L18").NumberFormat = "mmm-yy"
$ws.Range("M18").Value = 42887
$ws.Range("M18").NumberFormat = "mmm-yy"
$ws.Range("N18").Value = 27
$ws.Range("O18").Value = 1200
$ws.Range("P18").Value = 1.7

# --- Row 19 : Re-date ---
$ws.Range("A19").Value = "Re-date"
$ws.Range("L19").Value = 41518
$ws.Range("L19").NumberFormat = "mmm-yy"
$ws.Range("M19").Value = 42156
$ws.Range("M19").NumberFormat = "mmm-yy"
$ws.Range("N19").Value = 31
$ws.Range("O19").Value = 700
$ws.Range("P19").Value = 1.5

# Move the active selection to A20, matching where work continues next.
$ws.Range("A20").Select()
